$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3256
$ws.Range("B3").Value = 678
$ws.Range("B4").Value = 356
$ws.Range("B5").Value = 270
$ws.Range("B6").Value = 35
$ws.Range("B7").Value = 33
$ws.Range("B9").Value = 9
$ws.Range("B10").Value = 4

$ws.Range("A12").Value = "diorite gabbro anorthosite"
$ws.Range("A13").Value = "quartz syenite"
